$wb = $excel.ActiveWorkbook

$wsCreate = $wb.Worksheets.Item("Create")
$wsEdit   = $wb.Worksheets.Item("Edit")
$wsDelete = $wb.Worksheets.Item("Delete")

# --- Update cell values (order matters so new shared-strings are
#     allocated in the same sequence as the target workbook) ---

# "Edit" sheet: Updated Start Time changes from 08:50:00 -> 08:40:00
$wsEdit.Range("C2").Value = "08:40:00"

# Color Code changes from #61b377 -> #0e101e on all three sheets
$wsCreate.Range("C2").Value = "#0e101e"
$wsEdit.Range("D2").Value = "#0e101e"
$wsDelete.Range("C2").Value = "#0e101e"

# Updated Color Code changes from #304078 -> #1e0e16 on Edit/Delete sheets
$wsEdit.Range("E2").Value = "#1e0e16"
$wsDelete.Range("D2").Value = "#1e0e16"

# "Delete" sheet: Start Time changes from 08:50:00 -> 08:48:00
$wsDelete.Range("A2").Value = "08:48:00"

# "Create" sheet row 3: add a quote-prefixed time-formatted text value
$wsCreate.Range("A3").NumberFormat = "h:mm:ss"
$wsCreate.Range("A3").Value = "'08:00:00"

# --- Update selections on each sheet ---
$wsEdit.Range("A2").Select()
$wsDelete.Range("A2").Select()

# --- Switch the active tab from "Edit" to "Create", and select A3 there ---
$wsCreate.Activate()
$wsCreate.Range("A3").Select()
